# "Add files via upload" / "Update all maths of tables."
#
# The Futuros sheet has two example tables (a LONG table, rows 25-48, and a
# SHORT table, rows 52-75). This edit:
#   1. Clears the sample/example input numbers that were left in the
#      template (quantity, entry price, leverage, exit price - and for the
#      extra example rows, also the date/ticker columns), so the template
#      starts "empty" instead of pre-filled with worked examples.
#   2. Fixes the SL% (N column) formula so it divides the stop-loss % by
#      the leverage (instead of incorrectly treating $C$9 as a raw percent
#      literal).
#   3. Fixes the LONG table's "%"" (L) column formula so that it is
#      computed directly from the trade prices instead of back-solving
#      from the PNL-$ column.
#   4. For the SHORT table, the first row's "%" formula becomes a
#      standalone formula (no longer the head of the shared range), the
#      shared "%" range now starts on the following row, and the PNL-$
#      (M) column formula drops the extra leverage multiplication.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Futuros")

# ---------------------------------------------------------------------
# LONG table (rows 25:48)
# ---------------------------------------------------------------------

# Row 25 keeps its date/ticker/entry-time/exit-time (B,C,G,K) but loses the
# example quantity/price/leverage numbers.
$ws.Range("F25").ClearContents()
$ws.Range("H25:J25").ClearContents()

# Rows 26 and 27 lose every example input (date, ticker, quantity, times,
# price, leverage).
$ws.Range("B26:C27").ClearContents()
$ws.Range("F26:J26").ClearContents()
$ws.Range("G26").ClearContents()
$ws.Range("F27:J27").ClearContents()
$ws.Range("G27").ClearContents()
$ws.Range("K26:K27").ClearContents()

# Fixed "%" formula: was (M/E), now computed straight from the trade prices.
$ws.Range("L25:L48").Formula = "=((J25-H25)/H25)*I25"

# Fixed SL formula: was H*(1-$C$9%), now divides the stop-loss % by the
# leverage.
$ws.Range("N25:N48").Formula = "=H25*(1-`$C`$9/(I25*100))"

# ---------------------------------------------------------------------
# SHORT table (rows 52:75)
# ---------------------------------------------------------------------

# Row 52 keeps its entry/exit time (G,K) and ticker (C) but loses the
# example quantity/price/leverage numbers.
$ws.Range("F52").ClearContents()
$ws.Range("H52:J52").ClearContents()

# First row's "%" becomes a standalone (non-shared) formula.
$ws.Range("L52").Formula = "=((H52-J52)/H52)*I52"

# The shared "%" formula now starts on the next row instead.
$ws.Range("L53:L75").Formula = "=(M53/E53)"

# PNL-$ formula drops the extra "*I" (leverage) factor.
$ws.Range("M52:M75").Formula = "=(H52-J52)*F52"

# Fixed SL formula: was H*(1+$C$9%), now divides the stop-loss % by the
# leverage.
$ws.Range("N52:N75").Formula = "=H52*(1+(`$C`$9/(I52*100)))"
